$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.087.21'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.674.17'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.90'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5234'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.32%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2670'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06272'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.21'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07602'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.689.99'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.493'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5657'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008088'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.39'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.101.64'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.824'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.58'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '188.10'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.176'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.00%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.56'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1249'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.619'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.51%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06330'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.354'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.282'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.532'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.519'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.658'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.009'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.420'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6005'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.136'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01607'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.092.12'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8691'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.005'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.88'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.826.29'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000111'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.89'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.004'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05244'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.974'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.85%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.923'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.55%  '
